$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = -13.3508
$ws.Range("B9").Value = 6.084799999999994
$ws.Range("C12").Value = -11.133
$ws.Range("E13").Value = 16.66200000000001
$ws.Range("D15").Value = -8.835799999999994
$ws.Range("E16").Value = 16.44450000000001
$ws.Range("B18").Value = 6.960799999999998
$ws.Range("B20").Value = 9.595899999999993
$ws.Range("E20").Value = 15.95859999999999
$ws.Range("E24").Value = 16.68370000000001
$ws.Range("C26").Value = -12.5865
$ws.Range("B27").Value = 6.108500000000003
$ws.Range("C27").Value = -12.71219999999999
$ws.Range("C29").Value = -11.40600000000001
$ws.Range("C37").Value = -13.7998
$ws.Range("C38").Value = -13.0115
$ws.Range("D38").Value = -9.125899999999991
$ws.Range("E39").Value = 16.2818
$ws.Range("D44").Value = -7.387400000000001
$ws.Range("E48").Value = 17.39320000000002
$ws.Range("C51").Value = -12.0652
$ws.Range("D51").Value = -7.878200000000001
$ws.Range("E52").Value = 17.0223
$ws.Range("C55").Value = -13.92280000000001
$ws.Range("E56").Value = 16.56550000000001
$ws.Range("D57").Value = -8.358600000000001
$ws.Range("D63").Value = -7.4526
$ws.Range("B69").Value = 5.793100000000001
$ws.Range("C69").Value = -11.7188
$ws.Range("C70").Value = -12.2833
$ws.Range("D70").Value = -8.0436
$ws.Range("B76").Value = 5.2745
$ws.Range("B82").Value = 7.369900000000004
$ws.Range("C83").Value = -14.19839999999999
$ws.Range("E84").Value = 17.0015
$ws.Range("D99").Value = -8.134499999999997
$ws.Range("E100").Value = 16.383
$ws.Range("E101").Value = 16.77940000000001
$ws.Range("C102").Value = -13.16910000000001
